# Update "Forecast Comparison" sheet (Prophet/Amazon forecast numbers)
# These numbers changed after removing the Auto-ARIMA model from the forecast pipeline
$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$wsForecast.Range("C2").Value = 119
$wsForecast.Range("D2").Value = 105
$wsForecast.Range("E2").Value = 128
$wsForecast.Range("F2").Value = 155
$wsForecast.Range("G2").Value = 197
$wsForecast.Range("C3").Value = 123
$wsForecast.Range("D3").Value = 112
$wsForecast.Range("E3").Value = 136
$wsForecast.Range("F3").Value = 164
$wsForecast.Range("G3").Value = 209
$wsForecast.Range("C4").Value = 127
$wsForecast.Range("D4").Value = 124
$wsForecast.Range("E4").Value = 151
$wsForecast.Range("F4").Value = 182
$wsForecast.Range("G4").Value = 231
$wsForecast.Range("C5").Value = 123
$wsForecast.Range("D5").Value = 126
$wsForecast.Range("E5").Value = 153
$wsForecast.Range("F5").Value = 186
$wsForecast.Range("G5").Value = 239
$wsForecast.Range("C6").Value = 110
$wsForecast.Range("D6").Value = 131
$wsForecast.Range("E6").Value = 160
$wsForecast.Range("F6").Value = 196
$wsForecast.Range("G6").Value = 254
$wsForecast.Range("C7").Value = 107
$wsForecast.Range("D7").Value = 130
$wsForecast.Range("E7").Value = 159
$wsForecast.Range("F7").Value = 195
$wsForecast.Range("G7").Value = 252
$wsForecast.Range("C8").Value = 118
$wsForecast.Range("D8").Value = 129
$wsForecast.Range("E8").Value = 158
$wsForecast.Range("F8").Value = 194
$wsForecast.Range("G8").Value = 253
$wsForecast.Range("C9").Value = 130
$wsForecast.Range("D9").Value = 130
$wsForecast.Range("E9").Value = 159
$wsForecast.Range("F9").Value = 195
$wsForecast.Range("G9").Value = 253
$wsForecast.Range("C10").Value = 157
$wsForecast.Range("D10").Value = 126
$wsForecast.Range("E10").Value = 154
$wsForecast.Range("F10").Value = 189
$wsForecast.Range("G10").Value = 245
$wsForecast.Range("C11").Value = 180
$wsForecast.Range("D11").Value = 125
$wsForecast.Range("E11").Value = 153
$wsForecast.Range("F11").Value = 188
$wsForecast.Range("G11").Value = 245
$wsForecast.Range("C12").Value = 206
$wsForecast.Range("D12").Value = 123
$wsForecast.Range("E12").Value = 150
$wsForecast.Range("F12").Value = 186
$wsForecast.Range("G12").Value = 245
$wsForecast.Range("C13").Value = 219
$wsForecast.Range("D13").Value = 124
$wsForecast.Range("E13").Value = 152
$wsForecast.Range("F13").Value = 190
$wsForecast.Range("G13").Value = 250
$wsForecast.Range("C14").Value = 208
$wsForecast.Range("D14").Value = 122
$wsForecast.Range("E14").Value = 149
$wsForecast.Range("F14").Value = 185
$wsForecast.Range("G14").Value = 243
$wsForecast.Range("C15").Value = 169
$wsForecast.Range("D15").Value = 115
$wsForecast.Range("E15").Value = 141
$wsForecast.Range("F15").Value = 177
$wsForecast.Range("G15").Value = 236
$wsForecast.Range("C16").Value = 103
$wsForecast.Range("D16").Value = 116
$wsForecast.Range("E16").Value = 142
$wsForecast.Range("F16").Value = 177
$wsForecast.Range("G16").Value = 234
$wsForecast.Range("C17").Value = 59
$wsForecast.Range("D17").Value = 113
$wsForecast.Range("E17").Value = 139
$wsForecast.Range("F17").Value = 174
$wsForecast.Range("G17").Value = 231

# Update "Summary" sheet totals/extremes that are derived from the forecast values above
$wsSummary = $wb.Worksheets.Item("Summary")

# Keep these cells as text (they were stored as text in the original workbook)
$summaryCells = @("B9", "B10", "B11", "B12", "B14")
foreach ($addr in $summaryCells) {
    $wsSummary.Range($addr).NumberFormat = "@"
}

$wsSummary.Range("B9").Value  = "2258"   # Total Forecast (16 Weeks)
$wsSummary.Range("B10").Value = "957"    # Total Forecast (8 Weeks)
$wsSummary.Range("B11").Value = "492"    # Total Forecast (4 Weeks)
$wsSummary.Range("B12").Value = "219"    # Max Forecast
$wsSummary.Range("B14").Value = "59"     # Min Forecast

